# Add 2022-Q4 data
#
# This script:
#  1. Inserts a new worksheet named "2022-Q4" right after "总计", populating
#     it with the Q4 fund-holding detail rows (mirrors the layout used by the
#     other quarterly sheets).
#  2. Updates the "总计" (summary) sheet: inserts a new row for "2022-Q4" right
#     after the header, pushing the existing quarter rows down by one.
#  3. Restores the originally-active sheet ("2020-Q4", the last tab) so the
#     workbook/sheet view-selection state isn't disturbed by adding a sheet.

$wb = $excel.ActiveWorkbook
$summary = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# 1. Insert the new "2022-Q4" worksheet right after "总计".
# ---------------------------------------------------------------------------
$q3 = $wb.Worksheets.Item("2022-Q3")
$newSheet = $wb.Worksheets.Add($null, $summary)
$newSheet.Name = "2022-Q4"

# Copy the header-row formatting (bold / border / centered) from the existing
# "2022-Q3" sheet so the new sheet matches the look of its siblings.
$q3.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)

# Copy the "index" column formatting (A column, style used on A2 of "2022-Q3")
# across the 8 data rows we are about to fill in.
$q3.Range("A2").Copy()
$newSheet.Range("A2:A9").PasteSpecial(-4122)

# Header texts.
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Columns B..G hold fund codes / sizes / percentages that are stored as TEXT
# (not numbers) in this workbook (e.g. "000586", "2.30"). Force the range to
# text formatting first so the literal strings (incl. leading zeros / trailing
# zeros) are preserved instead of being auto-coerced into numbers.
$newSheet.Range("B2:G9").NumberFormat = "@"

$q4data = @(
    @("000586", "景顺长城中小创精选股票", "2.30", "89.60", "8.56", "0.1969", 1),
    @("260115", "景顺长城中小盘混合", "1.06", "91.71", "6.82", "0.0723", 3),
    @("009189", "华宝成长策略混合", "1.75", "88.29", "2.35", "0.0411", 10),
    @("001255", "长城改革红利灵活配置混合", "1.21", "81.29", "3.15", "0.0381", 6),
    @("010706", "景顺长城景骊成长混合", "0.68", "68.49", "5.61", "0.0381", 2),
    @("002681", "金鹰元和灵活配置混合A", "0.30", "81.19", "4.99", "0.0150", 7),
    @("002682", "金鹰元和灵活配置混合C", "0.23", "81.19", "4.99", "0.0115", 7),
    @("002303", "金鹰智慧生活灵活配置混合", "0.09", "94.65", "6.28", "0.0057", 2)
)

for ($i = 0; $i -lt $q4data.Length; $i++) {
    $row = $i + 2
    $rec = $q4data[$i]
    $newSheet.Cells.Item($row, 1).Value = $i          # column A: 0-based index
    $newSheet.Cells.Item($row, 2).Value = $rec[0]     # B: 基金代码
    $newSheet.Cells.Item($row, 3).Value = $rec[1]     # C: 基金名称
    $newSheet.Cells.Item($row, 4).Value = $rec[2]     # D: 基金规模
    $newSheet.Cells.Item($row, 5).Value = $rec[3]     # E: 股票总仓位
    $newSheet.Cells.Item($row, 6).Value = $rec[4]     # F: 仓位占比
    $newSheet.Cells.Item($row, 7).Value = $rec[5]     # G: 持有市值(亿元)
    $newSheet.Cells.Item($row, 8).Value = $rec[6]     # H: 仓位排名 (numeric)
}

# ---------------------------------------------------------------------------
# 2. Update the "总计" sheet: insert the 2022-Q4 row after the header, and
#    shift the existing rows down by one (same values/order as before).
# ---------------------------------------------------------------------------
$summaryRows = @(
    @("2022-Q4", 8, 0.42),
    @("2022-Q3", 2, 0.18),
    @("2022-Q2", 3, 0.24),
    @("2022-Q1", 3, 0.26),
    @("2021-Q4", 2, 0.26),
    @("2021-Q1", 3, 1.12),
    @("2020-Q4", 2, 1.01)
)

for ($i = 0; $i -lt $summaryRows.Length; $i++) {
    $row = $i + 2
    $rec = $summaryRows[$i]
    $summary.Cells.Item($row, 1).Value = $i      # column A: 0-based index
    $summary.Cells.Item($row, 2).Value = $rec[0] # column B: quarter label
    $summary.Cells.Item($row, 3).Value = $rec[1] # column C: holding count
    $summary.Cells.Item($row, 4).Value = $rec[2] # column D: holding value
}

# Row 8 is brand new (sheet used to stop at row 7) -- give its index cell (A8)
# the same style as the other index cells (A2:A7) by copying formats over.
$summary.Range("A7").Copy()
$summary.Range("A8").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 3. Restore original tab-selection state: "2020-Q4" (now the last sheet) was
#    the active tab before this edit; adding a worksheet switches the active
#    tab to the new one, so switch back.
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$lastSheet.Activate()

Write-Output "2022-Q4 sheet added; 总计 updated."
